# Insert a new data row above row 99 (this shifts existing rows 99-223
# down to 100-224, enlarging the used range from A1:R223 to A1:R224),
# then populate the newly inserted row with its data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(99).Insert()

$ws.Range("A99").Value = 5
$ws.Range("B99").Value = "Macroferia Regional de Talca"
$ws.Range("C99").Value = "Maule"
$ws.Range("D99").Value = 44671
$ws.Range("E99").Value = 7
$ws.Range("F99").Value = 100112008
$ws.Range("G99").Value = "Coliflor"
$ws.Range("H99").Value = "Sin especificar"
$ws.Range("I99").Value = "Primera"
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1000
$ws.Range("L99").Value = 1000
$ws.Range("M99").Value = 1000
$ws.Range("N99").Value = "$/unidad"
$ws.Range("O99").Value = "Región del Maule"
$ws.Range("P99").Value = 1000
$ws.Range("Q99").Value = 1
$ws.Range("R99").Value = "Hortaliza"
